$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. "Students" sheet: drop the ECTS_deficit column (E) - its data is no
#    longer tracked on this sheet. Clear contents only (keeps the column's
#    width metadata instead of shifting columns left).
$wsStudents = $wb.Worksheets.Item("Students")
$wsStudents.Activate()
$wsStudents.Range("E1:E31").ClearContents()
$wsStudents.Range("F4").Select()

# 2. "Grades comments" sheet: give every database row its own surrogate ID
#    column, consistent with the rest of the tables (class-per-table
#    refactor instead of a generic Table class).
$wsComments = $wb.Worksheets.Item("Grades comments")
$wsComments.Activate()
$wsComments.Columns.Item(1).Insert()
$wsComments.Range("A1").Value = "ID"
$wsComments.Range("A2").Value = 1
$wsComments.Columns.Item(1).ColumnWidth = 2.0
$wsComments.Columns.Item(2).ColumnWidth = 8.0
$wsComments.Columns.Item(3).ColumnWidth = 14.3
$wsComments.Range("D6").Select()

# 3. "Grades" sheet: no structural change, just where the cursor ended up
#    while reviewing the refactor.
$wsGrades = $wb.Worksheets.Item("Grades")
$wsGrades.Activate()
$wsGrades.Range("F1").Select()

# 4. Drop the "Suspended students" sheet entirely (its Reason/ECTS_deficit
#    fields are no longer modeled).
$wsSuspended = $wb.Worksheets.Item("Suspended students")
$wsSuspended.Delete()

# Leave "Students" as the active tab/selection, matching the saved file.
$wsStudents.Activate()
$wsStudents.Range("F4").Select()
